$wb = $excel.ActiveWorkbook

# ALC row 20
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(20, 8).Value = 26510.5  # H20: 80000 -> 26510.5
$ws.Cells.Item(20, 9).Value = 1021  # I20: 0 -> 1021
$ws.Cells.Item(20, 10).Value = 52000  # J20: 80000 -> 52000
$ws.Cells.Item(20, 11).Value = 1021  # K20: 0 -> 1021
$ws.Cells.Item(20, 12).Value = 52000  # L20: 80000 -> 52000
$ws.Cells.Item(20, 13).Value = -791  # M20: None -> -791
$ws.Cells.Item(20, 14).Value = -52460  # N20: -80460 -> -52460

# ALC row 35
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(35, 8).Value = 26510.5  # H35: 80000 -> 26510.5
$ws.Cells.Item(35, 9).Value = 1021  # I35: 0 -> 1021
$ws.Cells.Item(35, 10).Value = 52000  # J35: 80000 -> 52000
$ws.Cells.Item(35, 11).Value = 1021  # K35: 0 -> 1021
$ws.Cells.Item(35, 12).Value = 52000  # L35: 80000 -> 52000
$ws.Cells.Item(35, 13).Value = -642  # M35: None -> -642
$ws.Cells.Item(35, 14).Value = -52758  # N35: -80758 -> -52758

# ALC row 38
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(38, 8).Value = 649  # H38: 3766 -> 649
$ws.Cells.Item(38, 9).Value = 649  # I38: 649 -> 649
$ws.Cells.Item(38, 10).Value = 0  # J38: 10000 -> 0
$ws.Cells.Item(38, 11).Value = 1947  # K38: 1947 -> 1947
$ws.Cells.Item(38, 12).Value = 0  # L38: 30000 -> 0
$ws.Cells.Item(38, 13).Value = -1575  # M38: -1575 -> -1575
$ws.Cells.Item(38, 14).ClearContents()  # N38: -30744 -> (removed)

# ALC row 94
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(94, 8).Value = 1577.1111  # H94: 1677.4445 -> 1577.1111
$ws.Cells.Item(94, 9).Value = 1123.5  # I94: 1165.6666 -> 1123.5
$ws.Cells.Item(94, 10).Value = 1940  # J94: 1933.3334 -> 1940
$ws.Cells.Item(94, 11).Value = 1123.5  # K94: 1165.6666 -> 1123.5
$ws.Cells.Item(94, 12).Value = 1940  # L94: 1933.3334 -> 1940
$ws.Cells.Item(94, 13).Value = -672.5  # M94: -714.6666 -> -672.5
$ws.Cells.Item(94, 14).Value = -2842  # N94: -2835.3334 -> -2842

# ALC row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(100, 8).Value = 3087.5715  # H100: 2788.1538 -> 3087.5715
$ws.Cells.Item(100, 9).Value = 3466.2727  # I100: 3349.889 -> 3466.2727
$ws.Cells.Item(100, 10).Value = 1699  # J100: 1524.25 -> 1699
$ws.Cells.Item(100, 11).Value = 3466.2727  # K100: 3349.889 -> 3466.2727
$ws.Cells.Item(100, 12).Value = 1699  # L100: 1524.25 -> 1699
$ws.Cells.Item(100, 13).Value = -2925.2727  # M100: -2808.889 -> -2925.2727
$ws.Cells.Item(100, 14).Value = -2781  # N100: -2606.25 -> -2781

# ALC row 101
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(101, 8).Value = 5998  # H101: 0 -> 5998
$ws.Cells.Item(101, 9).Value = 5998  # I101: 0 -> 5998
$ws.Cells.Item(101, 10).Value = 0  # J101: 0 -> 0
$ws.Cells.Item(101, 11).Value = 17994  # K101: 0 -> 17994
$ws.Cells.Item(101, 12).Value = 0  # L101: 0 -> 0
$ws.Cells.Item(101, 13).Value = -16372  # M101: None -> -16372

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 1107  # H137: 883.4666999999999 -> 1107
$ws.Cells.Item(137, 9).Value = 845.2857  # I137: 674.9167 -> 845.2857
$ws.Cells.Item(137, 10).Value = 1717.6666  # J137: 1717.6666 -> 1717.6666
$ws.Cells.Item(137, 11).Value = 2535.8571  # K137: 2024.7501 -> 2535.8571
$ws.Cells.Item(137, 12).Value = 5152.9998  # L137: 5152.9998 -> 5152.9998
$ws.Cells.Item(137, 13).Value = 14.14289999999983  # M137: 525.2499 -> 14.14289999999983
$ws.Cells.Item(137, 14).Value = -10252.9998  # N137: -10252.9998 -> -10252.9998

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 20134.572  # H74: 24189.4 -> 20134.572
$ws.Cells.Item(74, 9).Value = 20134.572  # I74: 24189.4 -> 20134.572
$ws.Cells.Item(74, 10).Value = 0  # J74: 0 -> 0
$ws.Cells.Item(74, 11).Value = 20134.572  # K74: 24189.4 -> 20134.572
$ws.Cells.Item(74, 12).Value = 0  # L74: 0 -> 0
$ws.Cells.Item(74, 13).Value = -19260.572  # M74: -23315.4 -> -19260.572

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 20134.572  # H77: 24189.4 -> 20134.572
$ws.Cells.Item(77, 9).Value = 20134.572  # I77: 24189.4 -> 20134.572
$ws.Cells.Item(77, 10).Value = 0  # J77: 0 -> 0
$ws.Cells.Item(77, 11).Value = 100672.86  # K77: 120947 -> 100672.86
$ws.Cells.Item(77, 12).Value = 0  # L77: 0 -> 0
$ws.Cells.Item(77, 13).Value = -96304.86  # M77: -116579 -> -96304.86

# ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(102, 8).Value = 875.58826  # H102: 951.13336 -> 875.58826
$ws.Cells.Item(102, 9).Value = 910.875  # I102: 951.13336 -> 910.875
$ws.Cells.Item(102, 10).Value = 311  # J102: 0 -> 311
$ws.Cells.Item(102, 11).Value = 910.875  # K102: 951.13336 -> 910.875
$ws.Cells.Item(102, 12).Value = 311  # L102: 0 -> 311
$ws.Cells.Item(102, 13).Value = 711.125  # M102: 670.86664 -> 711.125
$ws.Cells.Item(102, 14).Value = -3555  # N102: None -> -3555

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(110, 8).Value = 1395  # H110: 3169.1667 -> 1395
$ws.Cells.Item(110, 9).Value = 1426.6666  # I110: 1173 -> 1426.6666
$ws.Cells.Item(110, 10).Value = 1300  # J110: 13150 -> 1300
$ws.Cells.Item(110, 11).Value = 1426.6666  # K110: 1173 -> 1426.6666
$ws.Cells.Item(110, 12).Value = 1300  # L110: 13150 -> 1300
$ws.Cells.Item(110, 13).Value = 618.3334  # M110: 872 -> 618.3334
$ws.Cells.Item(110, 14).Value = -5390  # N110: -17240 -> -5390

# ARM row 139
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(139, 8).Value = 70000  # H139: 85000 -> 70000
$ws.Cells.Item(139, 9).Value = 0  # I139: 0 -> 0
$ws.Cells.Item(139, 10).Value = 70000  # J139: 85000 -> 70000
$ws.Cells.Item(139, 11).Value = 0  # K139: 0 -> 0
$ws.Cells.Item(139, 12).Value = 70000  # L139: 85000 -> 70000
$ws.Cells.Item(139, 14).Value = -80280  # N139: -95280 -> -80280

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 6845.625  # H86: 7552.143 -> 6845.625
$ws.Cells.Item(86, 9).Value = 3316.5  # I86: 3788.6667 -> 3316.5
$ws.Cells.Item(86, 10).Value = 10374.75  # J86: 10374.75 -> 10374.75
$ws.Cells.Item(86, 11).Value = 3316.5  # K86: 3788.6667 -> 3316.5
$ws.Cells.Item(86, 12).Value = 10374.75  # L86: 10374.75 -> 10374.75
$ws.Cells.Item(86, 13).Value = -2193.5  # M86: -2665.6667 -> -2193.5
$ws.Cells.Item(86, 14).Value = -12620.75  # N86: -12620.75 -> -12620.75

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(89, 8).Value = 6845.625  # H89: 7552.143 -> 6845.625
$ws.Cells.Item(89, 9).Value = 3316.5  # I89: 3788.6667 -> 3316.5
$ws.Cells.Item(89, 10).Value = 10374.75  # J89: 10374.75 -> 10374.75
$ws.Cells.Item(89, 11).Value = 16582.5  # K89: 18943.3335 -> 16582.5
$ws.Cells.Item(89, 12).Value = 51873.75  # L89: 51873.75 -> 51873.75
$ws.Cells.Item(89, 13).Value = -10966.5  # M89: -13327.3335 -> -10966.5
$ws.Cells.Item(89, 14).Value = -63105.75  # N89: -63105.75 -> -63105.75

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 2560.6365  # H105: 2543.0833 -> 2560.6365
$ws.Cells.Item(105, 9).Value = 2566.9  # I105: 2566.9 -> 2566.9
$ws.Cells.Item(105, 10).Value = 2498  # J105: 2424 -> 2498
$ws.Cells.Item(105, 11).Value = 2566.9  # K105: 2566.9 -> 2566.9
$ws.Cells.Item(105, 12).Value = 2498  # L105: 2424 -> 2498
$ws.Cells.Item(105, 13).Value = -819.9000000000001  # M105: -819.9000000000001 -> -819.9000000000001
$ws.Cells.Item(105, 14).Value = -5992  # N105: -5918 -> -5992

# CRP row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 179.8  # H22: 300 -> 179.8
$ws.Cells.Item(22, 9).Value = 233.33333  # I22: 300 -> 233.33333
$ws.Cells.Item(22, 10).Value = 99.5  # J22: 0 -> 99.5
$ws.Cells.Item(22, 11).Value = 233.33333  # K22: 300 -> 233.33333
$ws.Cells.Item(22, 12).Value = 99.5  # L22: 0 -> 99.5
$ws.Cells.Item(22, 13).Value = 116.66667  # M22: 50 -> 116.66667
$ws.Cells.Item(22, 14).Value = -799.5  # N22: None -> -799.5

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2412.3333  # H31: 2495 -> 2412.3333
$ws.Cells.Item(31, 9).Value = 1999.6666  # I31: 2000 -> 1999.6666
$ws.Cells.Item(31, 10).Value = 2825  # J31: 2825 -> 2825
$ws.Cells.Item(31, 11).Value = 1999.6666  # K31: 2000 -> 1999.6666
$ws.Cells.Item(31, 12).Value = 2825  # L31: 2825 -> 2825
$ws.Cells.Item(31, 13).Value = -1704.6666  # M31: -1705 -> -1704.6666
$ws.Cells.Item(31, 14).Value = -3415  # N31: -3415 -> -3415

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 2412.3333  # H34: 2495 -> 2412.3333
$ws.Cells.Item(34, 9).Value = 1999.6666  # I34: 2000 -> 1999.6666
$ws.Cells.Item(34, 10).Value = 2825  # J34: 2825 -> 2825
$ws.Cells.Item(34, 11).Value = 1999.6666  # K34: 2000 -> 1999.6666
$ws.Cells.Item(34, 12).Value = 2825  # L34: 2825 -> 2825
$ws.Cells.Item(34, 13).Value = -1797.6666  # M34: -1798 -> -1797.6666
$ws.Cells.Item(34, 14).Value = -3229  # N34: -3229 -> -3229

# CRP row 39
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(39, 8).Value = 4754  # H39: 7557 -> 4754
$ws.Cells.Item(39, 9).Value = 3051  # I39: 0 -> 3051
$ws.Cells.Item(39, 10).Value = 6457  # J39: 7557 -> 6457
$ws.Cells.Item(39, 11).Value = 3051  # K39: 0 -> 3051
$ws.Cells.Item(39, 12).Value = 6457  # L39: 7557 -> 6457
$ws.Cells.Item(39, 13).Value = -2660  # M39: None -> -2660
$ws.Cells.Item(39, 14).Value = -7239  # N39: -8339 -> -7239

# CRP row 49
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(49, 8).Value = 4754  # H49: 7557 -> 4754
$ws.Cells.Item(49, 9).Value = 3051  # I49: 0 -> 3051
$ws.Cells.Item(49, 10).Value = 6457  # J49: 7557 -> 6457
$ws.Cells.Item(49, 11).Value = 3051  # K49: 0 -> 3051
$ws.Cells.Item(49, 12).Value = 6457  # L49: 7557 -> 6457
$ws.Cells.Item(49, 13).Value = -2869  # M49: None -> -2869
$ws.Cells.Item(49, 14).Value = -6821  # N49: -7921 -> -6821

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 6599.6665  # H58: 8700 -> 6599.6665
$ws.Cells.Item(58, 9).Value = 2399.5  # I58: 2400 -> 2399.5
$ws.Cells.Item(58, 10).Value = 15000  # J58: 15000 -> 15000
$ws.Cells.Item(58, 11).Value = 2399.5  # K58: 2400 -> 2399.5
$ws.Cells.Item(58, 12).Value = 15000  # L58: 15000 -> 15000
$ws.Cells.Item(58, 13).Value = -2196.5  # M58: -2197 -> -2196.5
$ws.Cells.Item(58, 14).Value = -15406  # N58: -15406 -> -15406

# CRP row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(105, 8).Value = 4658.1665  # H105: 4750 -> 4658.1665
$ws.Cells.Item(105, 9).Value = 4612.25  # I105: 4750 -> 4612.25
$ws.Cells.Item(105, 10).Value = 4750  # J105: 4750 -> 4750
$ws.Cells.Item(105, 11).Value = 4612.25  # K105: 4750 -> 4612.25
$ws.Cells.Item(105, 12).Value = 4750  # L105: 4750 -> 4750
$ws.Cells.Item(105, 13).Value = -2865.25  # M105: -3003 -> -2865.25
$ws.Cells.Item(105, 14).Value = -8244  # N105: -8244 -> -8244

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(122, 8).Value = 3861.625  # H122: 3777 -> 3861.625
$ws.Cells.Item(122, 9).Value = 4398.6  # I122: 3856.1428 -> 4398.6
$ws.Cells.Item(122, 10).Value = 2966.6667  # J122: 3500 -> 2966.6667
$ws.Cells.Item(122, 11).Value = 13195.8  # K122: 11568.4284 -> 13195.8
$ws.Cells.Item(122, 12).Value = 8900.000100000001  # L122: 10500 -> 8900.000100000001
$ws.Cells.Item(122, 13).Value = -10745.8  # M122: -9118.428400000001 -> -10745.8
$ws.Cells.Item(122, 14).Value = -13800.0001  # N122: -15400 -> -13800.0001

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 6599.6665  # H136: 8700 -> 6599.6665
$ws.Cells.Item(136, 9).Value = 2399.5  # I136: 2400 -> 2399.5
$ws.Cells.Item(136, 10).Value = 15000  # J136: 15000 -> 15000
$ws.Cells.Item(136, 11).Value = 7198.5  # K136: 7200 -> 7198.5
$ws.Cells.Item(136, 12).Value = 45000  # L136: 45000 -> 45000
$ws.Cells.Item(136, 13).Value = -4648.5  # M136: -4650 -> -4648.5
$ws.Cells.Item(136, 14).Value = -50100  # N136: -50100 -> -50100

# CRP row 141
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(141, 8).Value = 54098.57  # H141: 53129 -> 54098.57
$ws.Cells.Item(141, 9).Value = 19000  # I141: 19000 -> 19000
$ws.Cells.Item(141, 10).Value = 59948.332  # J141: 59954.8 -> 59948.332
$ws.Cells.Item(141, 11).Value = 19000  # K141: 19000 -> 19000
$ws.Cells.Item(141, 12).Value = 59948.332  # L141: 59954.8 -> 59948.332
$ws.Cells.Item(141, 13).Value = -13820  # M141: -13820 -> -13820
$ws.Cells.Item(141, 14).Value = -70308.33199999999  # N141: -70314.8 -> -70308.33199999999

# CUL row 3
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 3878.5715  # H3: 3447.5 -> 3878.5715
$ws.Cells.Item(3, 9).Value = 3878.5715  # I3: 3447.5 -> 3878.5715
$ws.Cells.Item(3, 10).Value = 0  # J3: 0 -> 0
$ws.Cells.Item(3, 11).Value = 11635.7145  # K3: 10342.5 -> 11635.7145
$ws.Cells.Item(3, 12).Value = 0  # L3: 0 -> 0
$ws.Cells.Item(3, 13).Value = -11523.7145  # M3: -10230.5 -> -11523.7145

# CUL row 137
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(137, 8).Value = 2257  # H137: 2673.75 -> 2257
$ws.Cells.Item(137, 9).Value = 2257  # I137: 2673.75 -> 2257
$ws.Cells.Item(137, 10).Value = 0  # J137: 0 -> 0
$ws.Cells.Item(137, 11).Value = 6771  # K137: 8021.25 -> 6771
$ws.Cells.Item(137, 12).Value = 0  # L137: 0 -> 0
$ws.Cells.Item(137, 13).Value = -1671  # M137: -2921.25 -> -1671

# GSM row 62
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(62, 8).Value = 0  # H62: 20000 -> 0
$ws.Cells.Item(62, 9).Value = 0  # I62: 20000 -> 0
$ws.Cells.Item(62, 10).Value = 0  # J62: 0 -> 0
$ws.Cells.Item(62, 11).Value = 0  # K62: 20000 -> 0
$ws.Cells.Item(62, 12).Value = 0  # L62: 0 -> 0
$ws.Cells.Item(62, 13).ClearContents()  # M62: -19314 -> (removed)

# GSM row 65
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(65, 8).Value = 0  # H65: 20000 -> 0
$ws.Cells.Item(65, 9).Value = 0  # I65: 20000 -> 0
$ws.Cells.Item(65, 10).Value = 0  # J65: 0 -> 0
$ws.Cells.Item(65, 11).Value = 0  # K65: 60000 -> 0
$ws.Cells.Item(65, 12).Value = 0  # L65: 0 -> 0
$ws.Cells.Item(65, 13).ClearContents()  # M65: -56568 -> (removed)

# GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 3949.5  # H97: 5999 -> 3949.5
$ws.Cells.Item(97, 9).Value = 5999  # I97: 5999 -> 5999
$ws.Cells.Item(97, 10).Value = 1900  # J97: 0 -> 1900
$ws.Cells.Item(97, 11).Value = 5999  # K97: 5999 -> 5999
$ws.Cells.Item(97, 12).Value = 1900  # L97: 0 -> 1900
$ws.Cells.Item(97, 13).Value = -5503  # M97: -5503 -> -5503
$ws.Cells.Item(97, 14).Value = -2892  # N97: None -> -2892

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 2874  # H122: 2915.889 -> 2874
$ws.Cells.Item(122, 9).Value = 2448.8  # I122: 2457.1667 -> 2448.8
$ws.Cells.Item(122, 10).Value = 5000  # J122: 3833.3333 -> 5000
$ws.Cells.Item(122, 11).Value = 7346.400000000001  # K122: 7371.500100000001 -> 7346.400000000001
$ws.Cells.Item(122, 12).Value = 15000  # L122: 11499.9999 -> 15000
$ws.Cells.Item(122, 13).Value = -4896.400000000001  # M122: -4921.500100000001 -> -4896.400000000001
$ws.Cells.Item(122, 14).Value = -19900  # N122: -16399.9999 -> -19900

# GSM row 136
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(136, 8).Value = 39110.75  # H136: 40088.4 -> 39110.75
$ws.Cells.Item(136, 9).Value = 0  # I136: 0 -> 0
$ws.Cells.Item(136, 10).Value = 39110.75  # J136: 40088.4 -> 39110.75
$ws.Cells.Item(136, 11).Value = 0  # K136: 0 -> 0
$ws.Cells.Item(136, 12).Value = 117332.25  # L136: 120265.2 -> 117332.25
$ws.Cells.Item(136, 14).Value = -122432.25  # N136: -125365.2 -> -122432.25

# LTW row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 2397.6  # H16: 1840.4286 -> 2397.6
$ws.Cells.Item(16, 9).Value = 2886  # I16: 2428 -> 2886
$ws.Cells.Item(16, 10).Value = 444  # J16: 371.5 -> 444
$ws.Cells.Item(16, 11).Value = 2886  # K16: 2428 -> 2886
$ws.Cells.Item(16, 12).Value = 444  # L16: 371.5 -> 444
$ws.Cells.Item(16, 13).Value = -2716  # M16: -2258 -> -2716
$ws.Cells.Item(16, 14).Value = -784  # N16: -711.5 -> -784

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 1916.8  # H93: 2121.375 -> 1916.8
$ws.Cells.Item(93, 9).Value = 1852.7142  # I93: 1978.5 -> 1852.7142
$ws.Cells.Item(93, 10).Value = 2066.3333  # J93: 2550 -> 2066.3333
$ws.Cells.Item(93, 11).Value = 1852.7142  # K93: 1978.5 -> 1852.7142
$ws.Cells.Item(93, 12).Value = 2066.3333  # L93: 2550 -> 2066.3333
$ws.Cells.Item(93, 13).Value = -604.7141999999999  # M93: -730.5 -> -604.7141999999999
$ws.Cells.Item(93, 14).Value = -4562.3333  # N93: -5046 -> -4562.3333

# WVR row 4
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(4, 8).Value = 60917.816  # H4: 220096 -> 60917.816
$ws.Cells.Item(4, 9).Value = 60917.816  # I4: 220096 -> 60917.816
$ws.Cells.Item(4, 10).Value = 0  # J4: 0 -> 0
$ws.Cells.Item(4, 11).Value = 60917.816  # K4: 220096 -> 60917.816
$ws.Cells.Item(4, 12).Value = 0  # L4: 0 -> 0
$ws.Cells.Item(4, 13).Value = -60804.816  # M4: -219983 -> -60804.816

# WVR row 23
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(23, 8).Value = 600  # H23: 0 -> 600
$ws.Cells.Item(23, 9).Value = 0  # I23: 0 -> 0
$ws.Cells.Item(23, 10).Value = 600  # J23: 0 -> 600
$ws.Cells.Item(23, 11).Value = 0  # K23: 0 -> 0
$ws.Cells.Item(23, 12).Value = 600  # L23: 0 -> 600
$ws.Cells.Item(23, 14).Value = -1058  # N23: None -> -1058
